$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B7").Value = 6.057
$ws.Range("A9").Value = -21.752
$ws.Range("B12").Value = 5.57
$ws.Range("B14").Value = 6.04
$ws.Range("A18").Value = -22.156
$ws.Range("A20").Value = -20.394
$ws.Range("B26").Value = 5.761
$ws.Range("A27").Value = -21.188
$ws.Range("B27").Value = 5.694999999999999
$ws.Range("B29").Value = 6.039
$ws.Range("A35").Value = -19.823
$ws.Range("B37").Value = 8.376000000000001
$ws.Range("B38").Value = 5.548
$ws.Range("B51").Value = 6.145
$ws.Range("B52").Value = 5.337
$ws.Range("B55").Value = 5.867
$ws.Range("A69").Value = -21.61
$ws.Range("B69").Value = 5.953
$ws.Range("B70").Value = 5.401999999999999
$ws.Range("A76").Value = -20.46
$ws.Range("A78").Value = -20.242
$ws.Range("B81").Value = 6.273
$ws.Range("A82").Value = -21.986
$ws.Range("A83").Value = -20.146
$ws.Range("B83").Value = 7.354000000000001
$ws.Range("A93").Value = -21.942
$ws.Range("B102").Value = 7.231
